$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=24; I="aa"; J="Agree/Accept"}
    @{Row=34; I="ba"; J="Appreciation"}
    @{Row=39; I="sv"; J="Statement-opinion"}
    @{Row=51; I="%"; J="Uninterpretable"}
    @{Row=52; I="aa"; J="Agree/Accept"}
    @{Row=57; I="sv"; J="Statement-opinion"}
    @{Row=62; I="ba"; J="Appreciation"}
    @{Row=69; I="sd"; J="Statement-non-opinion"}
    @{Row=75; I="aa"; J="Agree/Accept"}
    @{Row=76; I="%"; J="Uninterpretable"}
    @{Row=88; I="sd"; J="Statement-non-opinion"}
    @{Row=96; I="sv"; J="Statement-opinion"}
    @{Row=104; I="aa"; J="Agree/Accept"}
    @{Row=105; I="sd"; J="Statement-non-opinion"}
    @{Row=120; I="qy"; J="Yes-No-Question"}
    @{Row=125; I="aa"; J="Agree/Accept"}
    @{Row=129; I="sd"; J="Statement-non-opinion"}
    @{Row=162; I="sd"; J="Statement-non-opinion"}
    @{Row=173; I="aa"; J="Agree/Accept"}
    @{Row=183; I="sv"; J="Statement-opinion"}
    @{Row=184; I="sd"; J="Statement-non-opinion"}
    @{Row=186; I="sv"; J="Statement-opinion"}
    @{Row=190; I="sd"; J="Statement-non-opinion"}
    @{Row=191; I="sd"; J="Statement-non-opinion"}
    @{Row=196; I="ba"; J="Appreciation"}
    @{Row=197; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=205; I="ba"; J="Appreciation"}
    @{Row=206; I="sv"; J="Statement-opinion"}
    @{Row=207; I="ba"; J="Appreciation"}
    @{Row=209; I="sv"; J="Statement-opinion"}
    @{Row=220; I="sd"; J="Statement-non-opinion"}
    @{Row=223; I="aa"; J="Agree/Accept"}
    @{Row=239; I="aa"; J="Agree/Accept"}
    @{Row=241; I="sd"; J="Statement-non-opinion"}
    @{Row=242; I="sd"; J="Statement-non-opinion"}
    @{Row=251; I="ba"; J="Appreciation"}
    @{Row=267; I="%"; J="Uninterpretable"}
    @{Row=281; I="aa"; J="Agree/Accept"}
    @{Row=282; I="ba"; J="Appreciation"}
    @{Row=295; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=296; I="aa"; J="Agree/Accept"}
    @{Row=299; I="aa"; J="Agree/Accept"}
    @{Row=301; I="ba"; J="Appreciation"}
    @{Row=317; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=320; I="sv"; J="Statement-opinion"}
    @{Row=332; I="sv"; J="Statement-opinion"}
    @{Row=335; I="sd"; J="Statement-non-opinion"}
    @{Row=337; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=343; I="sv"; J="Statement-opinion"}
    @{Row=355; I="sv"; J="Statement-opinion"}
    @{Row=363; I="sv"; J="Statement-opinion"}
    @{Row=368; I="sd"; J="Statement-non-opinion"}
    @{Row=370; I="sd"; J="Statement-non-opinion"}
    @{Row=384; I="sv"; J="Statement-opinion"}
    @{Row=387; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=398; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=416; I="%"; J="Uninterpretable"}
    @{Row=432; I="sv"; J="Statement-opinion"}
    @{Row=441; I="sv"; J="Statement-opinion"}
    @{Row=447; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=448; I="sd"; J="Statement-non-opinion"}
    @{Row=453; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=460; I="sd"; J="Statement-non-opinion"}
    @{Row=479; I="%"; J="Uninterpretable"}
    @{Row=482; I="aa"; J="Agree/Accept"}
    @{Row=491; I="ba"; J="Appreciation"}
    @{Row=494; I="sd"; J="Statement-non-opinion"}
    @{Row=499; I="aa"; J="Agree/Accept"}
    @{Row=504; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=505; I="%"; J="Uninterpretable"}
    @{Row=512; I="sd"; J="Statement-non-opinion"}
    @{Row=527; I="sd"; J="Statement-non-opinion"}
    @{Row=531; I="ba"; J="Appreciation"}
    @{Row=532; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=533; I="aa"; J="Agree/Accept"}
    @{Row=537; I="ba"; J="Appreciation"}
    @{Row=538; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=542; I="sd"; J="Statement-non-opinion"}
    @{Row=547; I="sd"; J="Statement-non-opinion"}
    @{Row=553; I="ba"; J="Appreciation"}
    @{Row=554; I="sv"; J="Statement-opinion"}
    @{Row=564; I="sd"; J="Statement-non-opinion"}
    @{Row=567; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=576; I="%"; J="Uninterpretable"}
    @{Row=586; I="b"; J="Acknowledge (Backchannel)"}
    @{Row=600; I="aa"; J="Agree/Accept"}
    @{Row=610; I="aa"; J="Agree/Accept"}
    @{Row=619; I="ba"; J="Appreciation"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}

Write-Output "Updated $($updates.Count) rows"
